$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("G2").Style = "Normal"
$wsSummary.Range("D5").Select()

$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("N2").Copy($wsRepay.Range("O2"))
$wsRepay.Range("P2").Clear()
$wsRepay.Range("E12").Select()

$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 273
$wsTrans.Range("A3").Value = 272
$wsTrans.Range("A4").Value = 269
$wsTrans.Range("A5").Value = 271
$wsTrans.Range("A6").Value = 267
$wsTrans.Range("A7").Value = 266
$wsTrans.Range("D6").Select()
